$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 -> Sucrose__Ibark_Day_sp_exchange
$ws.Range("A2").Value = "Sucrose__Ibark_Day_sp_exchange"
$ws.Range("B2").Value = 0.01195219123505989
$ws.Range("C2").Value = -41.26647142419872
$ws.Range("D2").Value = 41.26647142419873
$ws.Range("E2").Value = $false

# Update row 3 -> L-Proline__Leaf_Day_sp_exchange
$ws.Range("A3").Value = "L-Proline__Leaf_Day_sp_exchange"
$ws.Range("B3").Value = 0.15
$ws.Range("C3").Value = -21.32001438946436
$ws.Range("D3").Value = 25.78744231144447
$ws.Range("E3").Value = $false

# Update row 4 -> L-Proline__Phellogen_Day_sp_exchange
$ws.Range("A4").Value = "L-Proline__Phellogen_Day_sp_exchange"
$ws.Range("B4").Value = 0.02777777777777782
$ws.Range("C4").Value = -12.53780877226151
$ws.Range("D4").Value = 12.56905456040514
$ws.Range("E4").Value = $false

# Update row 5 -> Citrate__Phellogen_Day_sp_exchange
$ws.Range("A5").Value = "Citrate__Phellogen_Day_sp_exchange"
$ws.Range("B5").Value = 0.02777777777777785
$ws.Range("C5").Value = -7.268536966848401
$ws.Range("D5").Value = 7.276804028092124
$ws.Range("E5").Value = $false

# Remove rows 6-13 which are no longer present in the data
$ws.Rows("6:13").Delete()
